$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Returns the Paragraph object containing the n-th (1-based) occurrence of
# $needle found in document order. Using text search (rather than a fixed
# Paragraphs(i) index) keeps this script resilient to the exact paragraph
# count of the document it is run against.
function Get-NthMatchParagraph($needle, $n) {
    $r = $d.Content.Duplicate
    $r.Start = 0
    $count = 0
    while ($r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $count += 1
        if ($count -eq $n) {
            return $r.Paragraphs(1)
        }
        $r.Start = $r.End
        $r.End = $d.Content.End
    }
    return $null
}

function Set-ParaXml($para, $innerXml) {
    $xml = '<w:p ' + $wns + '>' + $innerXml + '</w:p>'
    $para.Range.InsertXML($xml)
}

# --- DestinationFile="C:/MyFolder/My Report.rptdesign" />
# Drop the spellStart/spellEnd wrapping "MyFolder" and merge it with the
# following "/My " run into a single run.
$pMyFolder = Get-NthMatchParagraph 'MyFolder' 1
$xMyFolder = '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">                                     </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>DestinationFile</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>="C</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>:/</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve">MyFolder/My </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Report.rptdesign</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>" /&gt;</w:t></w:r>'
Set-ParaXml $pMyFolder $xMyFolder

# --- 1st "RolePermissions" - list item (R + olePermissions runs, wrapped in
# spellStart/spellEnd) -> "UserGroupPermissions" (UserGroup + Permissions
# runs, no proofErr at all).
$pRole1 = Get-NthMatchParagraph 'RolePermissions' 1
$xRole1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t>UserGroup</w:t></w:r>' +
    '<w:r><w:t>Permissions</w:t></w:r>'
Set-ParaXml $pRole1 $xRole1

# --- 2nd "RolePermissions" - attribute RolePermissions="Finance:VRE" (R +
# olePermissions runs) -> UserGroup + Permissions runs, keeping the
# spellStart/spellEnd wrap.
$pRole2 = Get-NthMatchParagraph 'RolePermissions' 1
$xRole2 = '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">                              </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>UserGroup</w:t></w:r>' +
    '<w:r><w:t>Permissions</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>="</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Finance</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>:VRE</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>"</w:t></w:r>'
Set-ParaXml $pRole2 $xRole2

# --- 3rd "RolePermissions" (now 1st again, since the two above were already
# rewritten and no longer match "RolePermissions") - list item single run
# "RolePermissions" -> split into UserGroup + Permissions runs, keeping the
# spellStart/spellEnd wrap and the trailing space run.
$pRole3 = Get-NthMatchParagraph 'RolePermissions' 1
$xRole3 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>UserGroup</w:t></w:r>' +
    '<w:r><w:t>Permissions</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParaXml $pRole3 $xRole3

# --- 4th "RolePermissions" (now the sole remaining match) - attribute
# RolePermissions="Finance:VRE" (R + olePermissions runs) -> UserGroup +
# Permissions runs, keeping the spellStart/spellEnd wrap.
$pRole4 = Get-NthMatchParagraph 'RolePermissions' 1
$xRole4 = '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">                             </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>UserGroup</w:t></w:r>' +
    '<w:r><w:t>Permissions</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>="</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Finance</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>:VRE</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>"</w:t></w:r>'
Set-ParaXml $pRole4 $xRole4

# --- TargetPath="C:/MyEnycDownloadDir" - drop the spellStart/spellEnd
# wrapping "MyEnycDownloadDir" and merge the three runs that spell out the
# value into one run.
$pEnyc = Get-NthMatchParagraph 'MyEnycDownloadDir' 1
$xEnyc = '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">                                             </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>TargetPath</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>="C:/MyEnycDownloadDir"</w:t></w:r>'
Set-ParaXml $pEnyc $xEnyc

# --- styles.xml: the TableNormal ("Normal Table") style picks up a
# <w:qFormat/> marker on resave.
$tableNormal = $d.Styles("Normal Table")
$tableNormal.QuickStyle = $true

Write-Output "done"
